{"js": "// Change \"\u00e9l\u00e9gant, l'application\" -> \"\u00e9l\u00e9gante, l'application\"\n// (agreement fix: \"\u00e9l\u00e9gant\" -> \"\u00e9l\u00e9gante\" before the feminine noun\n// \"l'application\"), leaving the rest of the paragraph untouched.\nconst body = context.document.body;\n\nconst searchText = \"\u00e9l\u00e9gant, l\\u2019application\";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found: \" + searchText);\n}\n\n// Replace the matched range's text, inserting the missing \"e\" so\n// \"\u00e9l\u00e9gant\" (masculine) becomes \"\u00e9l\u00e9gante\" (feminine), agreeing with\n// \"l'application\".\nconst target = results.items[0];\ntarget.insertText(\"\u00e9l\u00e9gante, l\\u2019application\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Change \"\u00e9l\u00e9gant, l'application\" -> \"\u00e9l\u00e9gante, l'application\"\n# (agreement fix: \"\u00e9l\u00e9gant\" -> \"\u00e9l\u00e9gante\" before the feminine noun\n# \"l'application\"), leaving the rest of the paragraph untouched.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u00e9l\u00e9gant, l\u2019application\"\n$find.Replacement.Text = \"\u00e9l\u00e9gante, l\u2019application\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchAllWordForms = $false\n\n# wdFindContinue = 1, wdReplaceOne = 1\n$found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n\nif (-not $found) {\n    throw \"Could not find target text '\u00e9l\u00e9gant, l\u2019application' to replace.\"\n}\n"}
